$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 16.62509999999999
$ws.Range("D7").Value = -6.989399999999997
$ws.Range("C8").Value = -12.5369
$ws.Range("C10").Value = -12.9301
$ws.Range("C12").Value = -10.83089999999999
$ws.Range("D15").Value = -8.109699999999995
$ws.Range("C18").Value = -14.0544
$ws.Range("D18").Value = -9.338199999999985
$ws.Range("E18").Value = 16.1819
$ws.Range("E19").Value = 16.5895
$ws.Range("D20").Value = -7.618799999999994
$ws.Range("E27").Value = 16.65659999999999
$ws.Range("D29").Value = -6.904300000000002
$ws.Range("D30").Value = -7.772999999999999
$ws.Range("D31").Value = -7.261999999999998
$ws.Range("E31").Value = 17.09890000000002
$ws.Range("C37").Value = -12.7639
$ws.Range("E38").Value = 16.23539999999999
$ws.Range("D40").Value = -8.520399999999992
$ws.Range("E42").Value = 16.4903
$ws.Range("E44").Value = 16.42809999999999
$ws.Range("E47").Value = 16.54429999999999
$ws.Range("D50").Value = -8.229399999999993
$ws.Range("C55").Value = -13.62769999999999
$ws.Range("E58").Value = 16.60390000000001
$ws.Range("E65").Value = 17.24650000000001
$ws.Range("C68").Value = -10.7614
$ws.Range("D68").Value = -6.998799999999997
$ws.Range("E73").Value = 17.27840000000002
$ws.Range("D76").Value = -7.554299999999999
$ws.Range("C77").Value = -12.2923
$ws.Range("C78").Value = -12.2676
$ws.Range("C81").Value = -12.9479
$ws.Range("C82").Value = -12.191
$ws.Range("D87").Value = -7.886099999999998
$ws.Range("D88").Value = -7.30139999999999
$ws.Range("E90").Value = 16.39759999999999
$ws.Range("E94").Value = 18.88720000000002
$ws.Range("E95").Value = 18.02630000000002
$ws.Range("D96").Value = -7.7723
$ws.Range("D98").Value = -8.532500000000004
$ws.Range("D101").Value = -7.737599999999997
$ws.Range("E101").Value = 16.5619
$ws.Range("D102").Value = -8.020699999999996
